$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 96
$ws.Range("I9").Value = 96.2
$ws.Range("J9").Value = 95
$ws.Range("K9").Value = 96.2
$ws.Range("L9").Value = 95
$ws.Range("M9").Value = 72.8
$ws.Range("N9").Value = -433

$ws.Range("H55").Value = 178.125
$ws.Range("I55").Value = 167.85715
$ws.Range("J55").Value = 250
$ws.Range("K55").Value = 167.85715
$ws.Range("L55").Value = 250
$ws.Range("M55").Value = 46.14285000000001
$ws.Range("N55").Value = -678

$ws.Range("H69").Value = 2980
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10748

$ws.Range("H72").Value = 2980
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35736

$ws.Range("H111").Value = 661.3333
$ws.Range("I111").Value = 393.8
$ws.Range("K111").Value = 1181.4
$ws.Range("M111").Value = 1885.6

$ws.Range("H132").Value = 911.2
$ws.Range("I132").Value = 873.6799999999999
$ws.Range("K132").Value = 2621.04
$ws.Range("M132").Value = -91.03999999999996

$ws.Range("H137").Value = 78466.08
$ws.Range("I137").Value = 1071
$ws.Range("J137").Value = 126838
$ws.Range("K137").Value = 3213
$ws.Range("L137").Value = 380514
$ws.Range("M137").Value = -663
$ws.Range("N137").Value = -385614

$ws.Range("H138").Value = 2097.3813
$ws.Range("I138").Value = 2252.2432
$ws.Range("J138").Value = 2001.8833
$ws.Range("K138").Value = 6756.7296
$ws.Range("L138").Value = 6005.6499
$ws.Range("M138").Value = -1616.7296
$ws.Range("N138").Value = -16285.6499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 618039.3
$ws.Range("I2").Value = 695044.25
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 695044.25
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -694931.25
$ws.Range("N2").Value = -2226

$ws.Range("H22").Value = 5001.067
$ws.Range("I22").Value = 5001.067
$ws.Range("K22").Value = 5001.067
$ws.Range("M22").Value = -4702.067

$ws.Range("H32").Value = 7624.3115
$ws.Range("I32").Value = 5870.3384
$ws.Range("J32").Value = 20876.555
$ws.Range("K32").Value = 5870.3384
$ws.Range("L32").Value = 20876.555
$ws.Range("M32").Value = -5583.3384
$ws.Range("N32").Value = -21450.555

$ws.Range("H74").Value = 2144.842
$ws.Range("I74").Value = 860.3
$ws.Range("J74").Value = 3572.111
$ws.Range("K74").Value = 860.3
$ws.Range("L74").Value = 3572.111
$ws.Range("M74").Value = 13.70000000000005
$ws.Range("N74").Value = -5320.111

$ws.Range("H77").Value = 2144.842
$ws.Range("I77").Value = 860.3
$ws.Range("J77").Value = 3572.111
$ws.Range("K77").Value = 4301.5
$ws.Range("L77").Value = 17860.555
$ws.Range("M77").Value = 66.5
$ws.Range("N77").Value = -26596.555

$ws.Range("H116").Value = 618039.3
$ws.Range("I116").Value = 695044.25
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 695044.25
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -692750.25
$ws.Range("N116").Value = -6588

$ws.Range("H122").Value = 1189.8948
$ws.Range("I122").Value = 1124
$ws.Range("K122").Value = 3372
$ws.Range("M122").Value = -922

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 618039.3
$ws.Range("I3").Value = 695044.25
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 695044.25
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -694930.25
$ws.Range("N3").Value = -2228

$ws.Range("H20").Value = 8579.429
$ws.Range("I20").Value = 2904.875
$ws.Range("J20").Value = 16145.5
$ws.Range("K20").Value = 2904.875
$ws.Range("L20").Value = 16145.5
$ws.Range("M20").Value = -2657.875
$ws.Range("N20").Value = -16639.5

$ws.Range("H105").Value = 2327.6155
$ws.Range("I105").Value = 2146.5
$ws.Range("K105").Value = 2146.5
$ws.Range("M105").Value = -399.5

$ws.Range("H134").Value = 6386.385
$ws.Range("I134").Value = 6787.478
$ws.Range("K134").Value = 20362.434
$ws.Range("M134").Value = -17827.434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 535.5
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H31").Value = 2256.4
$ws.Range("I31").Value = 1349.625
$ws.Range("J31").Value = 3292.7144
$ws.Range("K31").Value = 1349.625
$ws.Range("L31").Value = 3292.7144
$ws.Range("M31").Value = -1054.625
$ws.Range("N31").Value = -3882.7144

$ws.Range("H34").Value = 2256.4
$ws.Range("I34").Value = 1349.625
$ws.Range("J34").Value = 3292.7144
$ws.Range("K34").Value = 1349.625
$ws.Range("L34").Value = 3292.7144
$ws.Range("M34").Value = -1147.625
$ws.Range("N34").Value = -3696.7144

$ws.Range("H62").Value = 2902.5
$ws.Range("I62").Value = 2899
$ws.Range("K62").Value = 2899
$ws.Range("M62").Value = -2275

$ws.Range("H65").Value = 2902.5
$ws.Range("I65").Value = 2899
$ws.Range("K65").Value = 14495
$ws.Range("M65").Value = -11375

$ws.Range("H134").Value = 3499.8333
$ws.Range("I134").Value = 3499.8333
$ws.Range("K134").Value = 10499.4999
$ws.Range("M134").Value = -7964.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1064
$ws.Range("J107").Value = 1393.4166
$ws.Range("L107").Value = 4180.2498
$ws.Range("N107").Value = -8020.2498

$ws.Range("H122").Value = 678.0417
$ws.Range("I122").Value = 406
$ws.Range("K122").Value = 3654
$ws.Range("M122").Value = -1204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2651.9167
$ws.Range("I102").Value = 2772.3
$ws.Range("K102").Value = 2772.3
$ws.Range("M102").Value = -1150.3

$ws.Range("H107").Value = 163.25
$ws.Range("I107").Value = 150
$ws.Range("K107").Value = 150
$ws.Range("M107").Value = 1770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H122").Value = 61523.23
$ws.Range("I122").Value = 79584.8
$ws.Range("J122").Value = 1318
$ws.Range("K122").Value = 238754.4
$ws.Range("L122").Value = 3954
$ws.Range("M122").Value = -236304.4
$ws.Range("N122").Value = -8854
